$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.677041232585907
$ws.Range("B1").Value = 0.9734271764755249
$ws.Range("C1").Value = 4.142537117004395
$ws.Range("D1").Value = 2.173025846481323
$ws.Range("E1").Value = 1.657128810882568
